{"js": "// Insert a new paragraph \"Test paragraph\" right after the paragraph whose\n// text is \"Section 1\", using the same numbering (ilvl=0, numId=10) as that\n// paragraph, and explicit \"Normal\" paragraph style (pStyle).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the \"Section 1\" paragraph (anchor for the insertion).\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Section 1\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error('Could not find paragraph with text \"Section 1\"');\n}\n\n// Insert the new paragraph immediately after it.\nconst newPara = anchor.insertParagraph(\"Test paragraph\", Word.InsertLocation.after);\n\n// Explicitly apply the \"Normal\" paragraph style, then restore/apply the\n// list numbering (ilvl 0 / numId 10) that matches the anchor paragraph \u2014\n// setting the style resets any numbering picked up from the surrounding\n// context, so the list level must be (re)applied afterwards.\nnewPara.style = \"Normal\";\nnewPara.listItem.level = 0;\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"Test paragraph\" right after the paragraph whose\n# text is \"Section 1\", using the same numbering (ilvl=0, numId=10) as that\n# paragraph, and explicit \"Normal\" paragraph style (pStyle).\n\n$d = $word.ActiveDocument\n\n# Locate the \"Section 1\" paragraph (anchor for the insertion).\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text.TrimEnd(\"`r`a`v`n\") -eq \"Section 1\") {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not find paragraph with text 'Section 1'\"\n}\n\n$anchor = $d.Paragraphs($anchorIndex)\n$rng = $anchor.Range\n$rng.Collapse(0)  # wdCollapseEnd\n$rng.InsertAfter(\"Test paragraph`r\")\n\n# The newly inserted paragraph is the one right after the anchor.\n$newPara = $d.Paragraphs($anchorIndex + 1)\n\n# Applying the \"Normal\" style first (this clears any inherited numbering),\n# then (re)apply the list numbering level to match the anchor paragraph.\n$newPara.Style = \"Normal\"\n$newPara.Range.ListFormat.ListLevelNumber = 1\n"}
